# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (Home) target depth stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 221
$wsOff.Range("C2").Value = 168
$wsOff.Range("D2").Value = 70
$wsOff.Range("E2").Value = 37
$wsOff.Range("F2").Value = 4

# DEF sheet - row 2 (Home) target depth stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 228
$wsDef.Range("C2").Value = 175
$wsDef.Range("D2").Value = 66
$wsDef.Range("E2").Value = 24
